$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.399670124053955
$ws.Range("B1").Value = 4.149391651153564
$ws.Range("C1").Value = 0.1961207985877991
$ws.Range("D1").Value = 0.2006788700819016
$ws.Range("E1").Value = 0.3051038682460785
